$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "290.72"
Set-TextValue "E2" "-4.08%"
Set-TextValue "D3" "30.85"
Set-TextValue "E3" "-5.78%"
Set-TextValue "D4" "4.944"
Set-TextValue "E4" "0.37%"
Set-TextValue "D5" "0.07196"
Set-TextValue "E5" "-8.12%"
Set-TextValue "D6" "1.792"
Set-TextValue "E6" "-11.89%"
Set-TextValue "D7" "7.665"
Set-TextValue "E7" "-2.18%"
Set-TextValue "D8" "3.746"
Set-TextValue "E8" "-1.58%"
Set-TextValue "D9" "0.8955"
Set-TextValue "E9" "-2.91%"
Set-TextValue "D11" "0.07719"
Set-TextValue "E11" "-2.46%"
Set-TextValue "D12" "0.08071"
Set-TextValue "E12" "-6.64%"
Set-TextValue "D13" "0.03063"
Set-TextValue "E13" "-2.44%"
Set-TextValue "D14" "0.1003"
Set-TextValue "E14" "-0.22%"
Set-TextValue "D15" "0.001489"
Set-TextValue "E15" "-1.30%"
Set-TextValue "D16" "0.005769"
Set-TextValue "E16" "-2.59%"
Set-TextValue "D17" "3.472"
Set-TextValue "E17" "0.18%"
Set-TextValue "D18" "2.081"
Set-TextValue "E18" "-3.56%"
Set-TextValue "E19" "-0.86%"
Set-TextValue "E20" "-1.30%"
Set-TextValue "D21" "4.044"
Set-TextValue "E21" "-6.22%"
Set-TextValue "E22" "0.45%"
Set-TextValue "D23" "0.04523"
Set-TextValue "E23" "-1.13%"
Set-TextValue "D24" "0.001212"
Set-TextValue "E24" "-0.98%"
Set-TextValue "E25" "-9.91%"
Set-TextValue "D26" "0.0001250"
Set-TextValue "E26" "0.02%"
Set-TextValue "D39" "0.01603"
Set-TextValue "E39" "-7.80%"
Set-TextValue "D40" "0.04387"
Set-TextValue "E40" "-8.46%"
Set-TextValue "D41" "0.007355"
Set-TextValue "E41" "-2.11%"
Set-TextValue "D42" "0.1307"
Set-TextValue "E42" "-3.94%"
Set-TextValue "D43" "0.007671"
Set-TextValue "D44" "0.002040"
Set-TextValue "E44" "-12.80%"
Set-TextValue "D45" "0.009207"
Set-TextValue "E45" "-12.73%"
Set-TextValue "D46" "0.00005908"
Set-TextValue "E46" "-5.23%"
Set-TextValue "E47" "-0.01%"
Set-TextValue "E48" "173.66%"
Set-TextValue "D49" "0.003000"
Set-TextValue "E49" "-3.28%"
Set-TextValue "D50" "0.00002100"
Set-TextValue "E50" "-0.01%"
Set-TextValue "D51" "0.0002000"
Set-TextValue "E51" "-0.01%"
